$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content fix (commit: "case gender and published in objects corrected") ---
# B2 (CasesTab / "query" column): capitalization fix "Stage of Disease" -> "Stage Of Disease"
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis) 
 MATCH (samp:sample)-->(c)
 WHERE samp.specific_sample_pathology IN  ['T Cell Lymphoma']
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  
       coalesce(c.case_id, '') AS `Case ID`,
       coalesce(s.clinical_study_designation, '') AS `Study Code`,
       coalesce(s.clinical_study_type, '') AS  `Study Type`,
       coalesce(demo.breed, '') AS Breed ,
       coalesce(diag.disease_term, '') AS Diagnosis ,
       coalesce(diag.stage_of_disease, '') AS `Stage Of Disease`,
       CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
       coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
Order by c.case_id LIMIT 100
'@
$ws.Range("B2").Value = $casesQuery

# D2:D5 ("cartQuery" column, all tabs share the same cart query text): rewritten Cypher query
$cartQuery = @'
MATCH (samp:sample)-->(c:case)
WHERE 
  samp.specific_sample_pathology IN  ['T Cell Lymphoma']
WITH c
MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (f)-->(samp:sample)
OPTIONAL MATCH (f)-->(parent)
OPTIONAL MATCH (c)-->(s:study)
OPTIONAL MATCH (c)-->(cv:canine_individual)
WITH 
  DISTINCT f, samp, c, s, cv, parent,
  ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
  toInteger(floor(log(f.file_size)/log(1024))) as i,
  2 as precision
WITH 
  samp, c, s, f, cv, parent,
  f.file_size /(1024^i) AS value,
  10^precision AS factor,
  units[i] as unit
WITH
  samp, c, s, f, cv, parent, unit,
  round(factor * value)/factor AS size
RETURN
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_format, '') AS `Format`,
  coalesce(f.file_type, '') AS `File Type`,
  CASE size % 1 
      WHEN 0 
      THEN apoc.convert.toInteger(size)+' ' +unit 
      ELSE size+' ' +unit 
  END AS Size,
  head(labels(parent)) AS `Association`,
  coalesce(f.file_description,'') AS `Description`,
  coalesce(samp.sample_id, '') AS `Sample ID`,
  coalesce(c.case_id,'') as `Case ID`,
  coalesce(cv.canine_individual_id,'') AS `Canine ID`,
  CASE 
      WHEN s.clinical_study_designation IS NULL
      THEN parent.clinical_study_designation
      ELSE s.clinical_study_designation 
  END AS `Study Code`
ORDER BY `File Name` 
LIMIT 100
'@
$ws.Range("D2").Value = $cartQuery
$ws.Range("D3").Value = $cartQuery
$ws.Range("D4").Value = $cartQuery
$ws.Range("D5").Value = $cartQuery
